# This script applies the IFRS financial data corrections for JW Holdings
# (data/ifrs/JW홀딩스.xlsx): rows 2-6 (2014/12-2018/12, IFRS-consolidated
# actuals) get their reported figures replaced with the corrected values,
# and rows 7-9 (2019/12(E)-2021/12(E) consensus-estimate rows) have their
# numeric columns (D..AJ, excluding A/B/C) cleared, leaving only the
# row index / period-type / period-label columns populated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    "D2"=5739;
    "E2"=182;
    "F2"=182;
    "G2"=-131;
    "H2"=-93;
    "I2"=-109;
    "J2"=17;
    "K2"=10307;
    "L2"=8230;
    "M2"=2078;
    "N2"=939;
    "O2"=1138;
    "P2"=270;
    "Q2"=-172;
    "R2"=-19;
    "S2"=212;
    "T2"=96;
    "U2"=-268;
    "V2"=6648;
    "W2"=3.18;
    "X2"=-1.62;
    "Y2"=-10.98;
    "Z2"=-0.91;
    "AA2"=396.13;
    "AB2"=276.14;
    "AC2"=-177;
    "AD2"=-12.33;
    "AE2"=1530;
    "AF2"=1.42;
    "AG2"=22;
    "AH2"=1;
    "AI2"=-12.21;
    "AJ2"=61872025;
    "D3"=6228;
    "E3"=335;
    "F3"=335;
    "G3"=2;
    "H3"=3;
    "I3"=2;
    "J3"=2;
    "K3"=10829;
    "L3"=8126;
    "M3"=2703;
    "N3"=1037;
    "O3"=1666;
    "P3"=283;
    "Q3"=-240;
    "R3"=-382;
    "S3"=630;
    "T3"=77;
    "U3"=-316;
    "V3"=6486;
    "W3"=5.37;
    "X3"=0.05;
    "Y3"=0.18;
    "Z3"=0.03;
    "AA3"=300.61;
    "AB3"=288.37;
    "AC3"=3;
    "AD3"=2808.26;
    "AE3"=1685;
    "AF3"=4.89;
    "AG3"=41;
    "AH3"=0.5;
    "AI3"=1396.77;
    "AJ3"=61872025;
    "D4"=6848;
    "E4"=433;
    "F4"=433;
    "G4"=114;
    "H4"=-144;
    "I4"=-87;
    "J4"=-57;
    "K4"=11966;
    "L4"=9190;
    "M4"=2776;
    "N4"=1166;
    "O4"=1610;
    "P4"=307;
    "Q4"=-292;
    "R4"=-510;
    "S4"=1189;
    "T4"=107;
    "U4"=-399;
    "V4"=6689;
    "W4"=6.33;
    "X4"=-2.1;
    "Y4"=-7.87;
    "Z4"=-1.26;
    "AA4"=331.08;
    "AB4"=284.41;
    "AC4"=-137;
    "AD4"=-60.45;
    "AE4"=1782;
    "AF4"=4.66;
    "AG4"=70;
    "AH4"=0.84;
    "AI4"=-52.84;
    "AJ4"=65688800;
    "D5"=6812;
    "E5"=417;
    "F5"=417;
    "G5"=586;
    "H5"=569;
    "I5"=540;
    "J5"=28;
    "K5"=13473;
    "L5"=9766;
    "M5"=3706;
    "N5"=1661;
    "O5"=2071;
    "P5"=310;
    "Q5"=105;
    "R5"=146;
    "S5"=-335;
    "T5"=173;
    "U5"=-68;
    "V5"=7304;
    "W5"=6.12;
    "X5"=8.35;
    "Y5"=38.22;
    "Z5"=4.47;
    "AA5"=263.53;
    "AB5"=442.54;
    "AC5"=816;
    "AD5"=9.789999999999999;
    "AE5"=2505;
    "AF5"=3.19;
    "AG5"=70;
    "AH5"=0.88;
    "AI5"=8.59;
    "AJ5"=66452155;
    "D6"=7254;
    "E6"=333;
    "F6"=333;
    "G6"=4;
    "H6"=364;
    "I6"=232;
    "K6"=12756;
    "L6"=8746;
    "M6"=4010;
    "N6"=1818;
    "P6"=316;
    "Q6"=1958;
    "R6"=-913;
    "S6"=-734;
    "T6"=224;
    "U6"=1733;
    "V6"=5911;
    "W6"=4.59;
    "X6"=5.02;
    "Y6"=13.32;
    "Z6"=2.78;
    "AA6"=218.12;
    "AB6"=497.91;
    "AC6"=349;
    "AD6"=19.91;
    "AE6"=2742;
    "AF6"=2.53;
    "AG6"=81;
    "AH6"=1.17;
    "AI6"=23.16;
    "AJ6"=66452156;
}

foreach ($key in $newValues.Keys) {
    $ws.Range($key).Value = $newValues[$key]
}

$clearCells = @('D7','E7','G7','H7','I7','K7','L7','M7','N7','P7','Q7','R7','S7','T7','U7','W7','X7','Y7','Z7','AA7','AC7','AD7','AE7','AF7','AG7','AH7','AI7','D8','E8','G8','H8','I8','K8','L8','M8','N8','P8','Q8','R8','S8','T8','U8','W8','X8','Y8','Z8','AA8','AC8','AD8','AE8','AF8','AG8','AH8','AI8','D9','E9','G9','H9','I9','K9','L9','M9','N9','P9','Q9','R9','S9','T9','U9','W9','X9','Y9','Z9','AA9','AC9','AD9','AE9','AF9','AG9','AH9','AI9')

foreach ($addr in $clearCells) {
    $ws.Range($addr).ClearContents()
}
